$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: header + boolean flag ("Faltam dados para todos os Estados") ---
$ws.Range("F1").Value = "Faltam dados para todos os Estados"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- Force column B (Ano) to be stored as literal text so date-like strings are not reinterpreted as real dates ---
$ws.Range("B2:B31").NumberFormat = "@"

# Row 2: Brasil / 01/01/2015
$ws.Range("A2").Value = "Brasil"
$ws.Range("B2").Value = "01/01/2015"
$ws.Range("C2").Value = "Furto de veículo"
$ws.Range("D2").Value = 88.79143846282491
$ws.Range("E2").Formula = "=" + [char]34 + [char]34
$ws.Range("F2").Value = $true

# Row 3: Brasil / 01/01/2016
$ws.Range("A3").Value = "Brasil"
$ws.Range("B3").Value = "01/01/2016"
$ws.Range("C3").Value = "Furto de veículo"
$ws.Range("D3").Value = 95.93428856244273
$ws.Range("E3").Formula = "=" + [char]34 + [char]34
$ws.Range("F3").Value = $true

# Row 4: Brasil / 01/01/2017
$ws.Range("A4").Value = "Brasil"
$ws.Range("B4").Value = "01/01/2017"
$ws.Range("C4").Value = "Furto de veículo"
$ws.Range("D4").Value = 95.60128843421953
$ws.Range("E4").Formula = "=" + [char]34 + [char]34
$ws.Range("F4").Value = $true

# Row 5: Brasil / 01/01/2018
$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = "01/01/2018"
$ws.Range("C5").Value = "Furto de veículo"
$ws.Range("D5").Value = 92.59566656638793
$ws.Range("E5").Formula = "=" + [char]34 + [char]34
$ws.Range("F5").Value = $true

# Row 6: Brasil / 01/01/2019
$ws.Range("A6").Value = "Brasil"
$ws.Range("B6").Value = "01/01/2019"
$ws.Range("C6").Value = "Furto de veículo"
$ws.Range("D6").Value = 89.78982744962971
$ws.Range("E6").Formula = "=" + [char]34 + [char]34
$ws.Range("F6").Value = $true

# Row 7: Brasil / 01/01/2020
$ws.Range("A7").Value = "Brasil"
$ws.Range("B7").Value = "01/01/2020"
$ws.Range("C7").Value = "Furto de veículo"
$ws.Range("D7").Value = 71.10954397309881
$ws.Range("E7").Formula = "=" + [char]34 + [char]34
$ws.Range("F7").Value = $true

# Row 8: Brasil / 01/01/2021
$ws.Range("A8").Value = "Brasil"
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "Furto de veículo"
$ws.Range("D8").Value = 72.35437189743179
$ws.Range("E8").Formula = "=" + [char]34 + [char]34
$ws.Range("F8").Value = $false

# Row 9: Brasil / 01/01/2022
$ws.Range("A9").Value = "Brasil"
$ws.Range("B9").Value = "01/01/2022"
$ws.Range("C9").Value = "Furto de veículo"
$ws.Range("D9").Value = 81.75347312601198
$ws.Range("E9").Formula = "=" + [char]34 + [char]34
$ws.Range("F9").Value = $false

# Row 10: Brasil / 01/01/2023
$ws.Range("A10").Value = "Brasil"
$ws.Range("B10").Value = "01/01/2023"
$ws.Range("C10").Value = "Furto de veículo"
$ws.Range("D10").Value = 75.22949751208755
$ws.Range("E10").Formula = "=" + [char]34 + [char]34
$ws.Range("F10").Value = $false

# Row 11: Brasil / 01/01/2024
$ws.Range("A11").Value = "Brasil"
$ws.Range("B11").Value = "01/01/2024"
$ws.Range("C11").Value = "Furto de veículo"
$ws.Range("D11").Value = 70.4151726425975
$ws.Range("E11").Formula = "=" + [char]34 + [char]34
$ws.Range("F11").Value = $false

# Row 12: Nordeste / 01/01/2015
$ws.Range("A12").Value = "Nordeste"
$ws.Range("B12").Value = "01/01/2015"
$ws.Range("C12").Value = "Furto de veículo"
$ws.Range("D12").Value = 37.20037092649186
$ws.Range("E12").Formula = "=" + [char]34 + [char]34
$ws.Range("F12").Value = $true

# Row 13: Nordeste / 01/01/2016
$ws.Range("A13").Value = "Nordeste"
$ws.Range("B13").Value = "01/01/2016"
$ws.Range("C13").Value = "Furto de veículo"
$ws.Range("D13").Value = 47.47212669406449
$ws.Range("E13").Formula = "=" + [char]34 + [char]34
$ws.Range("F13").Value = $true

# Row 14: Nordeste / 01/01/2017
$ws.Range("A14").Value = "Nordeste"
$ws.Range("B14").Value = "01/01/2017"
$ws.Range("C14").Value = "Furto de veículo"
$ws.Range("D14").Value = 48.07072441239541
$ws.Range("E14").Formula = "=" + [char]34 + [char]34
$ws.Range("F14").Value = $true

# Row 15: Nordeste / 01/01/2018
$ws.Range("A15").Value = "Nordeste"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "Furto de veículo"
$ws.Range("D15").Value = 45.7067611166481
$ws.Range("E15").Formula = "=" + [char]34 + [char]34
$ws.Range("F15").Value = $true

# Row 16: Nordeste / 01/01/2019
$ws.Range("A16").Value = "Nordeste"
$ws.Range("B16").Value = "01/01/2019"
$ws.Range("C16").Value = "Furto de veículo"
$ws.Range("D16").Value = 45.2335334826059
$ws.Range("E16").Formula = "=" + [char]34 + [char]34
$ws.Range("F16").Value = $true

# Row 17: Nordeste / 01/01/2020
$ws.Range("A17").Value = "Nordeste"
$ws.Range("B17").Value = "01/01/2020"
$ws.Range("C17").Value = "Furto de veículo"
$ws.Range("D17").Value = 39.71933002184529
$ws.Range("E17").Formula = "=" + [char]34 + [char]34
$ws.Range("F17").Value = $true

# Row 18: Nordeste / 01/01/2021
$ws.Range("A18").Value = "Nordeste"
$ws.Range("B18").Value = "01/01/2021"
$ws.Range("C18").Value = "Furto de veículo"
$ws.Range("D18").Value = 39.17205746864656
$ws.Range("E18").Formula = "=" + [char]34 + [char]34
$ws.Range("F18").Value = $false

# Row 19: Nordeste / 01/01/2022
$ws.Range("A19").Value = "Nordeste"
$ws.Range("B19").Value = "01/01/2022"
$ws.Range("C19").Value = "Furto de veículo"
$ws.Range("D19").Value = 52.62425387813053
$ws.Range("E19").Formula = "=" + [char]34 + [char]34
$ws.Range("F19").Value = $false

# Row 20: Nordeste / 01/01/2023
$ws.Range("A20").Value = "Nordeste"
$ws.Range("B20").Value = "01/01/2023"
$ws.Range("C20").Value = "Furto de veículo"
$ws.Range("D20").Value = 54.12414104503122
$ws.Range("E20").Formula = "=" + [char]34 + [char]34
$ws.Range("F20").Value = $false

# Row 21: Nordeste / 01/01/2024
$ws.Range("A21").Value = "Nordeste"
$ws.Range("B21").Value = "01/01/2024"
$ws.Range("C21").Value = "Furto de veículo"
$ws.Range("D21").Value = 50.86320200298211
$ws.Range("E21").Formula = "=" + [char]34 + [char]34
$ws.Range("F21").Value = $false

# Row 22: Sergipe / 01/01/2015
$ws.Range("A22").Value = "Sergipe"
$ws.Range("B22").Value = "01/01/2015"
$ws.Range("C22").Value = "Furto de veículo"
$ws.Range("D22").Value = 36.2708348652949
$ws.Range("E22").Value = 18
$ws.Range("F22").Value = $true

# Row 23: Sergipe / 01/01/2016
$ws.Range("A23").Value = "Sergipe"
$ws.Range("B23").Value = "01/01/2016"
$ws.Range("C23").Value = "Furto de veículo"
$ws.Range("D23").Value = 50.10880001716483
$ws.Range("E23").Value = 18
$ws.Range("F23").Value = $true

# Row 24: Sergipe / 01/01/2017
$ws.Range("A24").Value = "Sergipe"
$ws.Range("B24").Value = "01/01/2017"
$ws.Range("C24").Value = "Furto de veículo"
$ws.Range("D24").Value = 38.01058448583375
$ws.Range("E24").Value = 22
$ws.Range("F24").Value = $true

# Row 25: Sergipe / 01/01/2018
$ws.Range("A25").Value = "Sergipe"
$ws.Range("B25").Value = "01/01/2018"
$ws.Range("C25").Value = "Furto de veículo"
$ws.Range("D25").Value = 29.05665081279616
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = $true

# Row 26: Sergipe / 01/01/2019
$ws.Range("A26").Value = "Sergipe"
$ws.Range("B26").Value = "01/01/2019"
$ws.Range("C26").Value = "Furto de veículo"
$ws.Range("D26").Value = 37.45266218394695
$ws.Range("E26").Value = 22
$ws.Range("F26").Value = $true

# Row 27: Sergipe / 01/01/2020
$ws.Range("A27").Value = "Sergipe"
$ws.Range("B27").Value = "01/01/2020"
$ws.Range("C27").Value = "Furto de veículo"
$ws.Range("D27").Value = 36.26513131341008
$ws.Range("E27").Value = 21
$ws.Range("F27").Value = $true

# Row 28: Sergipe / 01/01/2021
$ws.Range("A28").Value = "Sergipe"
$ws.Range("B28").Value = "01/01/2021"
$ws.Range("C28").Value = "Furto de veículo"
$ws.Range("D28").Value = 28.47750533632532
$ws.Range("E28").Value = 26
$ws.Range("F28").Value = $false

# Row 29: Sergipe / 01/01/2022
$ws.Range("A29").Value = "Sergipe"
$ws.Range("B29").Value = "01/01/2022"
$ws.Range("C29").Value = "Furto de veículo"
$ws.Range("D29").Value = 38.50982254945205
$ws.Range("E29").Value = 26
$ws.Range("F29").Value = $false

# Row 30: Sergipe / 01/01/2023
$ws.Range("A30").Value = "Sergipe"
$ws.Range("B30").Value = "01/01/2023"
$ws.Range("C30").Value = "Furto de veículo"
$ws.Range("D30").Value = 36.06223913262109
$ws.Range("E30").Value = 24
$ws.Range("F30").Value = $false

# Row 31: Sergipe / 01/01/2024
$ws.Range("A31").Value = "Sergipe"
$ws.Range("B31").Value = "01/01/2024"
$ws.Range("C31").Value = "Furto de veículo"
$ws.Range("D31").Value = 32.19900404429514
$ws.Range("E31").Value = 26
$ws.Range("F31").Value = $false
